$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "287.01"
$ws.Range("E2").Value = "2.14%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "29.45"
$ws.Range("E3").Value = "4.50%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.096"
$ws.Range("E4").Value = "1.18%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06950"
$ws.Range("E5").Value = "7.21%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "7.395"
$ws.Range("E6").Value = "2.17%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "3.537"
$ws.Range("E7").Value = "4.81%"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.15%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9024"
$ws.Range("E9").Value = "-3.07%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1587"
$ws.Range("E10").Value = "2.71%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06939"
$ws.Range("E11").Value = "13.15%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07613"
$ws.Range("E12").Value = "0.88%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02912"
$ws.Range("E13").Value = "1.47%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08984"
$ws.Range("E14").Value = "0.16%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001584"
$ws.Range("E15").Value = "0.19%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006475"
$ws.Range("E16").Value = "1.64%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006444"
$ws.Range("E17").Value = "6.58%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.469"
$ws.Range("E18").Value = "0.87%"

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "2.229"
$ws.Range("E19").Value = "-0.10%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3233"
$ws.Range("E20").Value = "1.35%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1324"
$ws.Range("E21").Value = "1.68%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "4.020"
$ws.Range("E22").Value = "-1.42%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.60%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04506"
$ws.Range("E24").Value = "1.51%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.87%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004377"
$ws.Range("E26").Value = "-0.35%"

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001166"
$ws.Range("E27").Value = "-6.91%"

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001610"
$ws.Range("E28").Value = "-0.42%"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.47%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006825"
$ws.Range("E41").Value = "2.64%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1243"
$ws.Range("E42").Value = "1.69%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002182"
$ws.Range("E43").Value = "7.81%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01156"
$ws.Range("E44").Value = "-4.26%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005742"
$ws.Range("E45").Value = "2.44%"

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01300"
$ws.Range("E47").Value = "0.11%"
